$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the custom (grey-fill) formatting on rows 28-29 so they match the
# default style used by the rest of the table.
$ws.Range("28:29").EntireRow.ClearFormats()

# Update the category code boundaries for the "Transfers" rows so the
# transfers adjustment is reflected in the monitoring table.
$ws.Range("D28").Value = 21050001
$ws.Range("C29").Value = 21050002

# Update the selected range to reflect the edited cells.
$ws.Range("C28:D29").Select()
